$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Price (D) column cells to Text format first so that
# numeric-looking strings (e.g. "1.00", "0.998") are preserved exactly
# as literal text instead of being coerced to numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D13","D14","D15","D16","D18","D19","D20","D22","D23","D24","D26","D27","D30","D31","D32","D33","D34","D35","D36","D39","D40","D41","D44","D45","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (price, volume%, and the two swapped rows coin/link/price cells).
$ws.Range("D2").Value = "57.681.54"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "2.329.07"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "541.74"
$ws.Range("E5").Value = "  +5.61%  "
$ws.Range("D6").Value = "134.89"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "2.360.00"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "2.764.52"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "23.52"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "57.696.79"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "2.348.12"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "10.57"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "334.85"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "6.75"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "61.79"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "8.44"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  +8.97%  "
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("D30").Value = "170.29"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "0.0₃0737"
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "18.57"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "1.02"
$ws.Range("E34").Value = "  +15.44%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +5.43%  "
$ws.Range("D39").Value = "1.61"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "39.45"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "149.98"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "284.31"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "19.30"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "17.59"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "0.380"
$ws.Range("E51").Value = "  -0.44%  "
